$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "25.992.05"
$ws.Range("E2").Value = "  +0.54%  "

Set-TextValue "D3" "1.741.81"
$ws.Range("E3").Value = "  +0.28%  "

Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue "D5" "246.38"
$ws.Range("E5").Value = "  +3.44%  "

$ws.Range("E6").Value = "  -0.04%  "

Set-TextValue "D7" "0.5023"
$ws.Range("E7").Value = "  -2.33%  "

Set-TextValue "D8" "0.2746"
$ws.Range("E8").Value = "  +1.00%  "

Set-TextValue "D9" "0.06187"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D10" "1.748.33"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D11" "0.07250"
$ws.Range("E11").Value = "  +1.17%  "

Set-TextValue "D12" "0.6533"
$ws.Range("E12").Value = "  +2.87%  "

$ws.Range("E13").Value = "  +1.15%  "

Set-TextValue "D14" "4.684"
$ws.Range("E14").Value = "  +2.09%  "

Set-TextValue "D15" "77.59"
$ws.Range("E15").Value = "  +0.71%  "

Set-TextValue "D16" "1.000"
$ws.Range("E16").Value = "  -0.16%  "

Set-TextValue "D17" "1.001"
$ws.Range("E17").Value = "  +0.02%  "

Set-TextValue "D18" "26.025.26"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  +1.82%  "

Set-TextValue "D20" "0.000006872"
$ws.Range("E20").Value = "  +2.75%  "

Set-TextValue "D21" "1.968.23"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  +5.17%  "

Set-TextValue "D23" "8.704"
$ws.Range("E23").Value = "  +0.86%  "

Set-TextValue "D24" "5.410"

Set-TextValue "D25" "135.59"
$ws.Range("E25").Value = "  -2.85%  "

Set-TextValue "D26" "1.512"
$ws.Range("E26").Value = "  +0.51%  "

Set-TextValue "D27" "15.26"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("E28").Value = "  +1.47%  "

Set-TextValue "D29" "105.64"
$ws.Range("E29").Value = "  +0.25%  "

Set-TextValue "D30" "3.948"
$ws.Range("E30").Value = "  +1.68%  "

Set-TextValue "D31" "0.08162"
$ws.Range("E31").Value = "  -2.13%  "

Set-TextValue "D32" "3.679"
$ws.Range("E32").Value = "  +2.87%  "

Set-TextValue "D33" "0.04696"
$ws.Range("E33").Value = "  +2.92%  "

Set-TextValue "D34" "2.667"
$ws.Range("E34").Value = "  +1.58%  "

Set-TextValue "D35" "0.9952"
$ws.Range("E35").Value = "  +1.21%  "

Set-TextValue "D36" "0.6111"
$ws.Range("E36").Value = "  -1.61%  "

Set-TextValue "D37" "2.762"
$ws.Range("E37").Value = "  +2.44%  "

Set-TextValue "D38" "0.01621"
$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  -0.03%  "

Set-TextValue "D41" "101.04"
$ws.Range("E41").Value = "  +3.62%  "

Set-TextValue "D42" "0.7939"
$ws.Range("E42").Value = "  +8.01%  "

Set-TextValue "D43" "0.3902"
$ws.Range("E43").Value = "  +1.47%  "

Set-TextValue "D44" "5.022"
$ws.Range("E44").Value = "  +1.83%  "

Set-TextValue "D45" "0.1165"
$ws.Range("E45").Value = "  +2.91%  "

Set-TextValue "D46" "6.329"
$ws.Range("E46").Value = "  +2.20%  "

Set-TextValue "D47" "55.65"
$ws.Range("E47").Value = "  +1.94%  "

Set-TextValue "D48" "0.05297"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "7.668"
$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D51" "0.3471"
$ws.Range("E51").Value = "  +1.65%  "
